$wb = $excel.ActiveWorkbook

# "Metadata" sheet: set the "Name" row's value (B4) to the ValueSet name
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B4").Value = "TypeidentifiantpersonneVs"

# Update the "Date" row's value (B8) to reflect the new generation timestamp
$wsMeta.Range("B8").Value = "2025-07-18T06:40:38+00:00"
